$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("kingdoms")

# Clear pre-existing formatting from these cells so they fall back to the
# default "Normal" style (no explicit style index), matching target output
$ws.Range("A4:C5").ClearFormats()

# Replicate the original authoring order so new shared-string indices line up
$ws.Range("C4").Value = "Lotharingia"
$ws.Range("B4").Value = "lotharingia"
$ws.Range("A4").Value = "(166, 235, 192)"
$ws.Range("A5").Value = "(116, 134, 194)"
$ws.Range("C5").Value = "Aquitaine"
$ws.Range("B5").Value = "aquitaine"
$ws.Range("E4").Value = "feudal_government"
$ws.Range("E5").Value = "feudal_government"

# Update selection to match diff
$ws.Range("B9").Select()
